$wb = $excel.ActiveWorkbook

# Sheet1 holds the pool of unused random names (one per row in column A).
# The name at row 2 ("j0writd2") has been consumed, so remove that row —
# every row below it shifts up by one and the sheet shrinks by one row.
$namesSheet = $wb.Worksheets.Item("Sheet1")
$namesSheet.Rows.Item(2).Delete()

# "used" sheet is an append-only log of consumed names. Record that
# "j0writd2" was just used for a newly generated image.
$usedSheet = $wb.Worksheets.Item("used")
$nextRow = $usedSheet.UsedRange.Rows.Count + 1

$usedSheet.Range("A$nextRow").Value = "j0writd2"
$usedSheet.Range("B$nextRow").Value = "ChatGPT Image 2026年1月20日 16_09_08.png"
$usedSheet.Range("C$nextRow").Value = "2026-01-20 16:10:24"
